$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DD")

# 1. Row 23: clear the "Trigger Creation Left" value from B23, keep its style (s=7)
$ws.Range("B23").ClearContents()

# 2. Row 28: add new "image_path" label in B28 (default style)
$ws.Range("B28").Value2 = "image_path"

# 3. Insert two new rows before the old row 39 ("Table: cart" block),
#    shifting all subsequent rows down by 2.
$ws.Rows("39:40").Insert()

$newCells = $ws.Range("B39:B40")
$newCells.Value2 = "productDescirption"
$ws.Range("B40").Value2 = "sub_category_id"

# Match the bordered look used elsewhere in column B (left+right thin border,
# same visual style as the neighboring "image_path" row, B38).
$newCells.Borders.Item(7).LineStyle = 1
$newCells.Borders.Item(7).Weight = 2
$newCells.Borders.Item(7).Color = 0
$newCells.Borders.Item(10).LineStyle = 1
$newCells.Borders.Item(10).Weight = 2
$newCells.Borders.Item(10).Color = 0

# 4. Restore the view/selection state recorded for this sheet after the edit.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("A44").Select()

Write-Host "Edit applied"
